$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15: Morning Glass of Ether
$ws.Range("H15").Value = 1250.23
$ws.Range("I15").Value = 1250.23
$ws.Range("K15").Value = 3750.69
$ws.Range("M15").Value = -3581.69

# ALC row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 874.3333
$ws.Range("I28").Value = 854.5
$ws.Range("K28").Value = 854.5
$ws.Range("M28").Value = -369.5

# ALC row 40: Stuck in the Moment
$ws.Range("H40").Value = 18519646
$ws.Range("I40").Value = 41667524
$ws.Range("K40").Value = 41667524
$ws.Range("M40").Value = -41667349

# ALC row 55: A Real Smooth Move
$ws.Range("H55").Value = 156.375
$ws.Range("I55").Value = 197.6
$ws.Range("J55").Value = 137.63637
$ws.Range("K55").Value = 197.6
$ws.Range("L55").Value = 137.63637
$ws.Range("M55").Value = 16.40000000000001
$ws.Range("N55").Value = -565.6363699999999

# ALC row 100: Asking for a Friend
$ws.Range("H100").Value = 1246.1538
$ws.Range("I100").Value = 1020
$ws.Range("K100").Value = 1020
$ws.Range("M100").Value = -479

# ALC row 116: Growing Up
$ws.Range("H116").Value = 2988.24
$ws.Range("I116").Value = 2400.3845
$ws.Range("K116").Value = 2400.3845
$ws.Range("M116").Value = 1041.6155

# ALC row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 336439.6
$ws.Range("I132").Value = 360321
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 1080963
$ws.Range("L132").Value = 6300
$ws.Range("M132").Value = -1078433
$ws.Range("N132").Value = -11360

# ALC row 141: Remedy for Reason
$ws.Range("H141").Value = 1758.6666
$ws.Range("I141").Value = 1588.0714
$ws.Range("J141").Value = 2099.8572
$ws.Range("K141").Value = 4764.2142
$ws.Range("L141").Value = 6299.571599999999
$ws.Range("M141").Value = 415.7857999999997
$ws.Range("N141").Value = -16659.5716

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45: Hollow Hallmarks
$ws.Range("H45").Value = 870
$ws.Range("I45").Value = 912.5
$ws.Range("J45").Value = 700
$ws.Range("K45").Value = 912.5
$ws.Range("L45").Value = 700
$ws.Range("M45").Value = -535.5
$ws.Range("N45").Value = -1454

# ARM row 138: Don't Ask about the Rivets
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
# BSM row 82: Spirituality Inspector
$ws.Range("H82").Value = 18194.77
$ws.Range("I82").Value = 10666.667
$ws.Range("J82").Value = 20453.2
$ws.Range("K82").Value = 10666.667
$ws.Range("L82").Value = 20453.2
$ws.Range("M82").Value = -10283.667
$ws.Range("N82").Value = -21219.2

# BSM row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 18194.77
$ws.Range("I85").Value = 10666.667
$ws.Range("J85").Value = 20453.2
$ws.Range("K85").Value = 10666.667
$ws.Range("L85").Value = 20453.2
$ws.Range("M85").Value = -9340.666999999999
$ws.Range("N85").Value = -23105.2

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16: Raise the Roof
$ws.Range("H16").Value = 2952.4211
$ws.Range("I16").Value = 3256.8572
$ws.Range("J16").Value = 2100
$ws.Range("K16").Value = 3256.8572
$ws.Range("L16").Value = 2100
$ws.Range("M16").Value = -2969.8572
$ws.Range("N16").Value = -2674

# CRP row 31: Wall Not Found
$ws.Range("H31").Value = 1338.68
$ws.Range("I31").Value = 1271.2106
$ws.Range("J31").Value = 1552.3334
$ws.Range("K31").Value = 1271.2106
$ws.Range("L31").Value = 1552.3334
$ws.Range("M31").Value = -976.2106000000001
$ws.Range("N31").Value = -2142.3334

# CRP row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1338.68
$ws.Range("I34").Value = 1271.2106
$ws.Range("J34").Value = 1552.3334
$ws.Range("K34").Value = 1271.2106
$ws.Range("L34").Value = 1552.3334
$ws.Range("M34").Value = -1069.2106
$ws.Range("N34").Value = -1956.3334

# CRP row 70: A Reward Fitting of the Faithful
$ws.Range("H70").Value = 28000
$ws.Range("J70").Value = 28000
$ws.Range("L70").Value = 28000
$ws.Range("N70").Value = -28630

# CRP row 73: Just Rewards for Just Devotion (L)
$ws.Range("H73").Value = 28000
$ws.Range("J73").Value = 28000
$ws.Range("L73").Value = 28000
$ws.Range("N73").Value = -30184

# CRP row 113: Patient Patients
$ws.Range("H113").Value = 2952.4211
$ws.Range("I113").Value = 3256.8572
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 3256.8572
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = -1086.8572
$ws.Range("N113").Value = -6440

$ws = $wb.Worksheets.Item("CUL")
# CUL row 34: Fever Pitch
$ws.Range("H34").Value = 62500880
$ws.Range("J34").Value = 66667596
$ws.Range("L34").Value = 200002788
$ws.Range("N34").Value = -200002956

# CUL row 88: Don't Let It Fall Apart
$ws.Range("H88").Value = 15000000
$ws.Range("I88").Value = 15000000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 45000000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -44999572
$ws.Range("N88").Value = ""

# CUL row 91: Better Come Back with a Sandwich (L)
$ws.Range("H91").Value = 15000000
$ws.Range("I91").Value = 15000000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 45000000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -44998518
$ws.Range("N91").Value = ""

# CUL row 113: Can't Eat Just One
$ws.Range("H113").Value = 588.619
$ws.Range("I113").Value = 583.8823
$ws.Range("J113").Value = 591.84
$ws.Range("K113").Value = 1751.6469
$ws.Range("L113").Value = 1775.52
$ws.Range("M113").Value = 418.3531
$ws.Range("N113").Value = -6115.52

# CUL row 117: A Good Omen
$ws.Range("H117").Value = 3971.4
$ws.Range("J117").Value = 6371.8335
$ws.Range("L117").Value = 19115.5005
$ws.Range("N117").Value = -25999.5005

# CUL row 129: Comfort Food
$ws.Range("H129").Value = 851.3889
$ws.Range("J129").Value = 1118.7273
$ws.Range("L129").Value = 3356.1819
$ws.Range("N129").Value = -13356.1819

# CUL row 131: The Mountain Steeped
$ws.Range("H131").Value = 2130570
$ws.Range("J131").Value = 2704935.2
$ws.Range("L131").Value = 8114805.600000001
$ws.Range("N131").Value = -8124885.600000001

# CUL row 133: Friends Are Food
$ws.Range("H133").Value = 4902.4287
$ws.Range("J133").Value = 4888.3887
$ws.Range("L133").Value = 14665.1661
$ws.Range("N133").Value = -24785.1661

# CUL row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 1237.0952
$ws.Range("I140").Value = 1251.9333
$ws.Range("K140").Value = 3755.7999
$ws.Range("M140").Value = 1424.2001

$ws = $wb.Worksheets.Item("GSM")
# GSM row 107: Whetstones for the Workers
$ws.Range("H107").Value = 612.1111
$ws.Range("I107").Value = 255.75
$ws.Range("J107").Value = 897.2
$ws.Range("K107").Value = 255.75
$ws.Range("L107").Value = 897.2
$ws.Range("M107").Value = 1664.25
$ws.Range("N107").Value = -4737.2

# GSM row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2091.923
$ws.Range("J113").Value = 2474.3333
$ws.Range("L113").Value = 2474.3333
$ws.Range("N113").Value = -6814.3333

# GSM row 126: Gold Rush Order
$ws.Range("H126").Value = 3290.4443
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 2722.8
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 8168.400000000001
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -13108.4

# GSM row 132: On Board for Lar
$ws.Range("H132").Value = 2079.8
$ws.Range("I132").Value = 1150
$ws.Range("J132").Value = 2699.6667
$ws.Range("K132").Value = 3450
$ws.Range("L132").Value = 8099.000100000001
$ws.Range("M132").Value = -920
$ws.Range("N132").Value = -13159.0001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7: Tan Before the Ban
$ws.Range("H7").Value = 2509.4
$ws.Range("I7").Value = 2232.6667
$ws.Range("K7").Value = 2232.6667
$ws.Range("M7").Value = -2120.6667

# LTW row 40: Best Served Toad
$ws.Range("H40").Value = 2831.25
$ws.Range("I40").Value = 2770.5715
$ws.Range("J40").Value = 2916.2
$ws.Range("K40").Value = 2770.5715
$ws.Range("L40").Value = 2916.2
$ws.Range("M40").Value = -2634.5715
$ws.Range("N40").Value = -3188.2

# LTW row 61: Spelling Me Softly
$ws.Range("H61").Value = 3179
$ws.Range("I61").Value = 2950
$ws.Range("J61").Value = 3331.6667
$ws.Range("K61").Value = 2950
$ws.Range("L61").Value = 3331.6667
$ws.Range("M61").Value = -2748
$ws.Range("N61").Value = -3735.6667

# LTW row 113: Peace in Rest
$ws.Range("H113").Value = 3179
$ws.Range("I113").Value = 2950
$ws.Range("J113").Value = 3331.6667
$ws.Range("K113").Value = 2950
$ws.Range("L113").Value = 3331.6667
$ws.Range("M113").Value = -780
$ws.Range("N113").Value = -7671.6667

# LTW row 126: Battered Books
$ws.Range("H126").Value = 2509.4
$ws.Range("I126").Value = 2232.6667
$ws.Range("K126").Value = 6698.000100000001
$ws.Range("M126").Value = -4228.000100000001

# LTW row 132: Tenets of Tanning
$ws.Range("H132").Value = 6851.6665
$ws.Range("I132").Value = 11173.454
$ws.Range("J132").Value = 3880.4375
$ws.Range("K132").Value = 33520.362
$ws.Range("L132").Value = 11641.3125
$ws.Range("M132").Value = -30990.362
$ws.Range("N132").Value = -16701.3125

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 5398
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 6122.5
$ws.Range("K81").Value = 5000
$ws.Range("L81").Value = 12245
$ws.Range("M81").Value = -3939
$ws.Range("N81").Value = -14367

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 5398
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 6122.5
$ws.Range("K84").Value = 25000
$ws.Range("L84").Value = 61225
$ws.Range("M84").Value = -19696
$ws.Range("N84").Value = -71833

# WVR row 113: A Tender Table
$ws.Range("H113").Value = 450.21738
$ws.Range("I113").Value = 458.35715
$ws.Range("J113").Value = 437.55554
$ws.Range("K113").Value = 1375.07145
$ws.Range("L113").Value = 1312.66662
$ws.Range("M113").Value = 794.9285500000001
$ws.Range("N113").Value = -5652.66662

# WVR row 132: Comfy Cabins
$ws.Range("H132").Value = 4226.409
$ws.Range("I132").Value = 5227.2856
$ws.Range("J132").Value = 2474.875
$ws.Range("K132").Value = 15681.8568
$ws.Range("L132").Value = 7424.625
$ws.Range("M132").Value = -13151.8568
$ws.Range("N132").Value = -12484.625

# WVR row 136: Weaving the Envelope
$ws.Range("H136").Value = 1262.6957
$ws.Range("I136").Value = 1192.4762
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3577.4286
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1027.4286
$ws.Range("N136").Value = -11100
